# resetSignal request: append a new day's attendance/SpO2 log sheet
# ("2020-11-08"), mirroring the existing "2020-11-02" sheet's layout.

$wb = $excel.ActiveWorkbook

# Find the last existing sheet so the new tab is appended at the end
# (after "2020-11-02"), matching the target sheet order/ids.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2020-11-08"

# Seed layout + formatting (header bold/border/alignment, styled column A)
# by copying the template sheet's first three rows, then overwrite the
# cells whose values actually differ for this day.
$template = $wb.Worksheets.Item("2020-11-02")
$template.Range("A1:F3").Copy($ws.Range("A1:F3"))

$ws.Range("E2").Value = "14:24:41"
$ws.Range("F2").Value = 93.10171335021791

$ws.Range("E3").Value = "14:26:27"
$ws.Range("F3").Value = 93.83166958817134
